$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove fully-deleted rows (from bottom to top to keep row numbers stable) ---
# Row 8 = "Dead adventurer" (冒险者尸体) - card removed
$ws.Range("A8:F8").EntireRow.Delete()
# Row 4 = "Corner" (拐角) - card removed
$ws.Range("A4:F4").EntireRow.Delete()
# Row 3 = "Hiding spot" (隐蔽处) - card removed
$ws.Range("A3:F3").EntireRow.Delete()

# At this point the rows are:
#  2 Exit/房间出口   3 杂物堆   4 Lever/拉杆   5 Chest/宝箱   6 Resource merchant/商人   7 Trainer/训练师   8 Altar/祭坛

# --- Step 2: reposition the Altar (祭坛) row ahead of the merchant/trainer rows ---
# Its content is being fully rewritten anyway, so drop the old row and open a fresh slot at row 6.
$ws.Range("A8:F8").EntireRow.Delete()
$ws.Range("A6:F6").EntireRow.Insert()

# Now the rows are:
#  2 Exit/房间出口   3 杂物堆   4 Lever/拉杆   5 Chest/宝箱   6 (blank)   7 Resource merchant/商人   8 Trainer/训练师

# --- Step 3: update text in rows that only need a rename / description edit ---
# Row 2: "Exit" card renamed to "Dungeon reset"; keep D2/E2 (cardNameEn/effectEn) untouched.
$ws.Range("A2").Value = '地下城重整'
$ws.Range("C2").Value = '重抽场上所有敌人侧的牌。'

# Row 3: junk pile description updated only.
$ws.Range("C3").Value = '翻开遭遇牌堆前3张牌，获得其中1张物品牌。'

# Row 4: lever description updated only; keep D4/E4 (Lever / choose-one rich text) untouched.
$ws.Range("C4").Value = '选场上1张陷阱牌，触发它或将其移动到场上任意位置。'

# Row 5: chest description updated only; keep D5/E5 (Chest english text) untouched.
$ws.Range("C5").Value = '获得遭遇牌堆第1张物品牌。'

# --- Step 4: fill the fresh Altar row (6) ---
# The freshly inserted row copied formatting from its neighbour, so start by
# wiping D6/E6 completely (no cardNameEn/effectEn for this card).
$ws.Range("D6:E6").Clear()
$ws.Range("A6").Value = '祭坛'
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = '多选：①弃置1张物品牌。②将1张手牌送墓。③受到2伤害。<br>
执行2项以上时：从购买能力区选1张牌获得。'
$ws.Range("C6").WrapText = $true

# --- Step 5: merchant -> shop rename + new text; maxCount 1 -> 2. Keep D7/E7 untouched. ---
$ws.Range("A7").Value = '商店'
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = '多选：①可重复，弃置1张物品牌，然后获得1金币。②取遭遇牌堆前3张物品牌，然后玩家每支付2金币，可以获得其中1张牌。'

# --- Step 6: trainer -> training ground rename + new text; maxCount 1 -> 2; drop cardNameEn/effectEn. ---
$ws.Range("A8").Value = '训练场'
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = '可重复：支付3金币，从购买能力区选1张牌获得。'
$ws.Range("D8:E8").Clear()

# --- Step 7: row heights follow Excel's normal auto-fit behaviour ---
$ws.Range("A1:F8").EntireRow.AutoFit()

# --- Step 8: sheet view bookkeeping to match the saved file ---
$ws.Range("D7").Select()
